$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J40").Value = 2244.4443
$ws.Range("I40").Value = 1789.8334
$ws.Range("K40").Value = 1789.8334
$ws.Range("M40").Value = -1614.8334
$ws.Range("L40").Value = 2244.4443
$ws.Range("H40").Value = 2062.6
$ws.Range("N40").Value = -2594.4443

$ws.Range("N103").Value = -4371.0002
$ws.Range("J103").Value = 1066.3334
$ws.Range("L103").Value = 3199.0002
$ws.Range("H103").Value = 100839.8

$ws.Range("N116").Value = -8526.5
$ws.Range("J116").Value = 1642.5
$ws.Range("L116").Value = 1642.5
$ws.Range("K116").Value = 5557473.5
$ws.Range("I116").Value = 5557473.5
$ws.Range("M116").Value = -5554031.5
$ws.Range("H116").Value = 5153413

$ws.Range("N132").Value = -20060.273
$ws.Range("J132").Value = 5000.091
$ws.Range("K132").Value = 12531.2901
$ws.Range("I132").Value = 4177.0967
$ws.Range("L132").Value = 15000.273
$ws.Range("M132").Value = -10001.2901
$ws.Range("H132").Value = 4392.643

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N61").Value = -2674
$ws.Range("J61").Value = 2250
$ws.Range("K61").Value = 1087
$ws.Range("I61").Value = 1087
$ws.Range("L61").Value = 2250
$ws.Range("M61").Value = -875
$ws.Range("H61").Value = 1232.375

$ws.Range("N74").Value = -4970
$ws.Range("K74").Value = 12503067
$ws.Range("J74").Value = 3222
$ws.Range("I74").Value = 12503067
$ws.Range("L74").Value = 3222
$ws.Range("M74").Value = -12502193
$ws.Range("H74").Value = 8067638

$ws.Range("H77").Value = 8067638
$ws.Range("N77").Value = -24846
$ws.Range("J77").Value = 3222
$ws.Range("K77").Value = 62515335
$ws.Range("I77").Value = 12503067
$ws.Range("L77").Value = 16110
$ws.Range("M77").Value = -62510967

$ws.Range("H122").Value = 3034.5
$ws.Range("N122").Value = -14894.2
$ws.Range("J122").Value = 3331.4
$ws.Range("I122").Value = 1550
$ws.Range("K122").Value = 4650
$ws.Range("L122").Value = 9994.200000000001
$ws.Range("M122").Value = -2200

$ws.Range("N132").Value = -12649.4
$ws.Range("J132").Value = 2529.8
$ws.Range("K132").Value = 4474.227000000001
$ws.Range("I132").Value = 1491.409
$ws.Range("L132").Value = 7589.400000000001
$ws.Range("M132").Value = -1944.227000000001
$ws.Range("H132").Value = 1815.9062

$ws.Range("J136").Value = 2250
$ws.Range("K136").Value = 3261
$ws.Range("L136").Value = 6750
$ws.Range("M136").Value = -711
$ws.Range("I136").Value = 1087
$ws.Range("H136").Value = 1232.375
$ws.Range("N136").Value = -11850

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("N58").Value = -1725.091
$ws.Range("J58").Value = 1319.091
$ws.Range("K58").Value = 583.36365
$ws.Range("I58").Value = 583.36365
$ws.Range("L58").Value = 1319.091
$ws.Range("M58").Value = -380.36365
$ws.Range("H58").Value = 877.65454

$ws.Range("J136").Value = 1319.091
$ws.Range("K136").Value = 1750.09095
$ws.Range("L136").Value = 3957.273
$ws.Range("M136").Value = 799.90905
$ws.Range("I136").Value = 583.36365
$ws.Range("H136").Value = 877.65454
$ws.Range("N136").Value = -9057.272999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N26").Value = -5560
$ws.Range("J26").Value = 5000
$ws.Range("L26").Value = 5000
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H26").Value = 5000

$ws.Range("L50").Value = 5000
$ws.Range("H50").Value = 5000
$ws.Range("I50").Value = 0
$ws.Range("N50").Value = -5996
$ws.Range("K50").Value = 0
$ws.Range("J50").Value = 5000
$ws.Range("M50").ClearContents()

$ws.Range("N132").Value = -16065.5
$ws.Range("J132").Value = 3668.5
$ws.Range("K132").Value = 4126.875
$ws.Range("I132").Value = 1375.625
$ws.Range("L132").Value = 11005.5
$ws.Range("M132").Value = -1596.875
$ws.Range("H132").Value = 1948.8438

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3334152.2
$ws.Range("M22").Value = -4762338
$ws.Range("N22").Value = -1620
$ws.Range("I22").Value = 4762633
$ws.Range("K22").Value = 4762633
$ws.Range("J22").Value = 1030
$ws.Range("L22").Value = 1030

$ws.Range("I27").Value = 4762633
$ws.Range("L27").Value = 1030
$ws.Range("M27").Value = -4762526
$ws.Range("H27").Value = 3334152.2
$ws.Range("N27").Value = -1244
$ws.Range("J27").Value = 1030
$ws.Range("K27").Value = 4762633

$ws.Range("J46").Value = 908.7895
$ws.Range("K46").Value = 763
$ws.Range("I46").Value = 763
$ws.Range("L46").Value = 908.7895
$ws.Range("M46").Value = -575
$ws.Range("H46").Value = 888.9091
$ws.Range("N46").Value = -1284.7895

$ws.Range("H128").Value = 30000
$ws.Range("N128").Value = -39960
$ws.Range("J128").Value = 30000
$ws.Range("L128").Value = 30000

$ws.Range("N132").Value = -11153.8574
$ws.Range("J132").Value = 2031.2858
$ws.Range("K132").Value = 39910386
$ws.Range("I132").Value = 13303462
$ws.Range("L132").Value = 6093.857400000001
$ws.Range("M132").Value = -39907856
$ws.Range("H132").Value = 11579203

$ws.Range("J136").Value = 1311.8462
$ws.Range("K136").Value = 21018.45
$ws.Range("L136").Value = 3935.5386
$ws.Range("M136").Value = -18468.45
$ws.Range("I136").Value = 7006.15
$ws.Range("H136").Value = 4762.9395
$ws.Range("N136").Value = -9035.5386

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1811.4445
$ws.Range("I122").Value = 1445.6666
$ws.Range("K122").Value = 4336.9998
$ws.Range("M122").Value = -1886.9998

$ws.Range("N132").Value = -14120.2499
$ws.Range("J132").Value = 3020.0833
$ws.Range("K132").Value = 2373.7143
$ws.Range("I132").Value = 791.2381
$ws.Range("L132").Value = 9060.249899999999
$ws.Range("M132").Value = 156.2856999999999
$ws.Range("H132").Value = 1601.7273

$ws.Range("J136").Value = 2045
$ws.Range("K136").Value = 1925.9031
$ws.Range("L136").Value = 6135
$ws.Range("M136").Value = 624.0969
$ws.Range("I136").Value = 641.9677
$ws.Range("H136").Value = 1138.875
$ws.Range("N136").Value = -11235
